$d = $word.ActiveDocument

$d.Content.Find.Execute("#NOME_CLIENTE", $true, $false, $false, $false, $false,
                         $true, 1, $false, "#NOME_CLIENTE", 2)
